# cv3 deck: drop the "enum" slide (it became redundant after the new
# "Ukoly"/"Ukol 3" pair was finalized). Everything else keeps its place;
# PowerPoint just shifts the remaining slides up by one and renumbers the
# relationship ids internally when it re-saves the package.

$p = $ppt.ActivePresentation

# Locate the "enum" slide (Title 1 == "enum") defensively instead of
# hard-coding index 7, in case slides ever get reordered upstream.
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = ""
    try { $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text } catch {}
    if ($title -eq "enum") {
        $target = $slide
        break
    }
}

if ($target -ne $null) {
    $target.Delete()
}
